$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: the original ORD_DET_01 test case, shifted down from row 2 ---
$ws.Range("A4").Value = "ORD_DET_01"
$ws.Range("B4").Value = "Xem chi tiết đơn hàng"
$ws.Range("C4").Value = "orderId = 100"
$ws.Range("D4").Value = "1. orderId=100`n2. Gọi service getDetails`n3. Forward JSP"
$ws.Range("E4").Value = "Forward → View-order-detail.jsp"
$ws.Range("F4").Value = "OK"
$ws.Range("G4").Value = "PASS"
$ws.Range("G4").Font.Bold = $true
$ws.Range("G4").Font.Name = "Calibri"
$ws.Range("G4").Font.Size = 11
$ws.Range("G4").Font.Color = 32768
$ws.Rows.Item(4).EntireRow.AutoFit()

# --- Row 2: new ORD_DET_03 test case (Service Exception) ---
$ws.Range("A2").Value = "ORD_DET_03"
$ws.Range("B2").Value = "Lỗi hệ thống (Service)"
$ws.Range("C2").Value = "Service Exception"
$ws.Range("D2").Value = "1. orderId=10`n2. Service ném lỗi`n3. Catch & Forward Error"
$ws.Range("E2").Value = "Forward → error.jsp"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"
$ws.Rows.Item(2).EntireRow.AutoFit()

# --- Row 3: new ORD_DET_02 test case (parse error) ---
$ws.Range("A3").Value = "ORD_DET_02"
$ws.Range("B3").Value = "ID lỗi (Chữ)"
$ws.Range("C3").Value = "orderId = abc"
$ws.Range("D3").Value = "1. orderId='abc'`n2. ParseInt lỗi`n3. Catch & Forward Error"
$ws.Range("E3").Value = "Forward → error.jsp"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"
$ws.Range("G3").Font.Bold = $true
$ws.Range("G3").Font.Name = "Calibri"
$ws.Range("G3").Font.Size = 11
$ws.Range("G3").Font.Color = 32768
$ws.Rows.Item(3).EntireRow.AutoFit()

# --- Column width tweaks (Dữ Liệu Mẫu / Các Bước columns) ---
$ws.Columns.Item(3).ColumnWidth = 15.8
$ws.Columns.Item(4).ColumnWidth = 22.5
